$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.413.86'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '1.823.01'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.13'
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4488'
$ws.Range("E7").Value = '  +2.11%  '
$ws.Range("E8").Value = '  +2.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07504'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8868'
$ws.Range("E10").Value = '  +5.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.06'
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '1.826.78'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.765'
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.03'
$ws.Range("E14").Value = '  +4.96%  '
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07103'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008797'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("D21").Value = '27.423.68'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.284'
$ws.Range("E22").Value = '  +2.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.93'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '2.055.54'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.955'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.369'
$ws.Range("E26").Value = '  +7.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.56'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.376'
$ws.Range("E29").Value = '  +2.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.09'
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08862'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7865'
$ws.Range("E32").Value = '  +6.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.520'
$ws.Range("E34").Value = '  +2.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.922'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9995'
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.113'
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05325'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.364'
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5327'
$ws.Range("E41").Value = '  +3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1724'
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.855'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.299'
$ws.Range("E44").Value = '  +19.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.741'
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5093'
$ws.Range("E46").Value = '  +5.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.68'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.76'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9996'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06382'
$ws.Range("E51").Value = '  +0.72%  '
